$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Rows 2-33: price (D) and volume-change (E) updates
Set-TextValue 2 4 "60.494.23"
Set-TextValue 2 5 "  -2.16%  "
Set-TextValue 3 4 "2.957.90"
Set-TextValue 3 5 "  -1.70%  "
Set-TextValue 4 4 "0.998"
Set-TextValue 4 5 "  -0.34%  "
Set-TextValue 5 4 "517.00"
Set-TextValue 5 5 "  -1.76%  "
Set-TextValue 6 4 "129.70"
Set-TextValue 6 5 "  +1.92%  "
Set-TextValue 7 4 "0.999"
Set-TextValue 7 5 "  -0.03%  "
Set-TextValue 8 4 "2.957.70"
Set-TextValue 8 5 "  -1.57%  "
Set-TextValue 9 5 "  -2.14%  "
Set-TextValue 10 4 "6.14"
Set-TextValue 10 5 "  +3.47%  "
Set-TextValue 11 4 "0.147"
Set-TextValue 11 5 "  -0.07%  "
Set-TextValue 12 4 "0.436"
Set-TextValue 12 5 "  -1.30%  "
Set-TextValue 13 5 "  -0.36%  "
Set-TextValue 14 4 "32.88"
Set-TextValue 14 5 "  -0.14%  "
Set-TextValue 15 4 "3.429.10"
Set-TextValue 15 5 "  -1.23%  "
Set-TextValue 16 4 "0.110"
Set-TextValue 16 5 "  +0.43%  "
Set-TextValue 17 4 "60.312.42"
Set-TextValue 17 5 "  -2.47%  "
Set-TextValue 18 4 "2.933.62"
Set-TextValue 18 5 "  -2.88%  "
Set-TextValue 19 4 "6.45"
Set-TextValue 19 5 "  +0.63%  "
Set-TextValue 20 4 "456.12"
Set-TextValue 20 5 "  -2.43%  "
Set-TextValue 21 4 "12.96"
Set-TextValue 21 5 "  +0.35%  "
Set-TextValue 22 4 "0.667"
Set-TextValue 22 5 "  -1.52%  "
Set-TextValue 23 4 "6.75"
Set-TextValue 23 5 "  -1.26%  "
Set-TextValue 24 4 "77.68"
Set-TextValue 24 5 "  +0.30%  "
Set-TextValue 25 4 "11.62"
Set-TextValue 25 5 "  +0.24%  "
Set-TextValue 26 5 "  +0.54%  "
Set-TextValue 27 4 "2.59"
Set-TextValue 27 5 "  -0.61%  "
Set-TextValue 28 4 "7.61"
Set-TextValue 28 5 "  -2.84%  "
Set-TextValue 29 5 "  -0.41%  "
Set-TextValue 30 4 "25.09"
Set-TextValue 30 5 "  -0.25%  "
Set-TextValue 31 5 "  +4.60%  "
Set-TextValue 32 4 "1.82"
Set-TextValue 32 5 "  +1.35%  "
Set-TextValue 33 4 "54.78"
Set-TextValue 33 5 "  -1.59%  "

# Rows 34 and 35: Stacks/NEARProtocol swap positions with updated data
Set-TextValue 34 2 "NEARProtocol"
Set-TextValue 34 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 34 4 "5.34"
Set-TextValue 34 5 "  +6.04%  "
Set-TextValue 35 2 "Stacks"
Set-TextValue 35 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue 35 4 "2.23"
Set-TextValue 35 5 "  -2.86%  "

# Rows 36-50: price (D) and volume-change (E) updates
Set-TextValue 36 4 "5.72"
Set-TextValue 36 5 "  -0.28%  "
Set-TextValue 37 4 "446.67"
Set-TextValue 37 5 "  -3.21%  "
Set-TextValue 38 4 "3.125.88"
Set-TextValue 38 5 "  +3.70%  "
Set-TextValue 39 4 "0.0772"
Set-TextValue 39 5 "  +0.49%  "
Set-TextValue 40 4 "0.0375"
Set-TextValue 40 5 "  -1.98%  "
Set-TextValue 41 4 "0.116"
Set-TextValue 41 5 "  +4.69%  "
Set-TextValue 42 4 "7.93"
Set-TextValue 42 5 "  +1.41%  "
Set-TextValue 43 4 "2.41"
Set-TextValue 43 5 "  -1.30%  "
Set-TextValue 44 5 "  +0.11%  "
Set-TextValue 45 5 "  -0.03%  "
Set-TextValue 46 4 "25.04"
Set-TextValue 46 5 "  +7.22%  "
Set-TextValue 47 4 "118.65"
Set-TextValue 47 5 "  +1.93%  "
Set-TextValue 48 4 "0.107"
Set-TextValue 48 5 "  +1.07%  "
Set-TextValue 49 4 "1.93"
Set-TextValue 49 5 "  -1.01%  "
Set-TextValue 50 4 "0.0₃0502"
Set-TextValue 50 5 "  -1.32%  "

# Row 51: BitgetToken replaced by CoreDAO with new data
Set-TextValue 51 2 "CoreDAO"
Set-TextValue 51 3 "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue 51 4 "2.26"
Set-TextValue 51 5 "  +0.27%  "
